# Update countries & provincias Spain
# - Reorder a handful of country-name rows (label swaps) so the country
#   list reflects the corrected alphabetical-ish placement.
# - Refresh COVID case numbers for several countries (data update).
# - Bump the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country label swaps (A column) -----------------------------------
$ws.Range("A139").Value = "Ruanda"
$ws.Range("A140").Value = "Principado de Andorra"

$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("A201").Value = "Laos"

$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"

$ws.Range("A212").Value = "Seychelles"
$ws.Range("A213").Value = "Montserrat"

# --- Updated case numbers ----------------------------------------------
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2530499
$ws.Range("C4").Value = 25911
$ws.Range("D4").Value = 1055328
$ws.Range("E4").Value = 1348052
$ws.Range("G4").Value = 339
$ws.Range("H4").Value = 127119

# Row 17 - Alemania
$ws.Range("B17").Value = 194042
$ws.Range("C17").Value = 257
$ws.Range("E17").Value = 7925
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 9017

# Row 19 - Francia
$ws.Range("B19").Value = 162936
$ws.Range("C19").Value = 1588
$ws.Range("E19").Value = 57807
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 29778

# Row 30 - Ecuador
$ws.Range("B30").Value = 53856
$ws.Range("C30").Value = 700
$ws.Range("D30").Value = 26493
$ws.Range("E30").Value = 22957
$ws.Range("G30").Value = 63
$ws.Range("H30").Value = 4406

# Row 50 - Barein
$ws.Range("E50").Value = 5507
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 73

# Row 68 - Marruecos
$ws.Range("B68").Value = 11633
$ws.Range("C68").Value = 295
$ws.Range("D68").Value = 8656
$ws.Range("E68").Value = 2759
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 218

# Row 75 - Uzbekistan
$ws.Range("B75").Value = 7320
$ws.Range("C75").Value = 143
$ws.Range("D75").Value = 5038
$ws.Range("E75").Value = 2262

# Row 107 - Mali
$ws.Range("B107").Value = 2060
$ws.Range("C107").Value = 21
$ws.Range("D107").Value = 1387
$ws.Range("E107").Value = 560

# Row 139 (now "Ruanda" after the label swap above)
$ws.Range("B139").Value = 858
$ws.Range("C139").Value = 8
$ws.Range("D139").Value = 398
$ws.Range("E139").Value = 458
$ws.Range("H139").Value = 2

# Row 140 (now "Principado de Andorra" after the label swap above)
$ws.Range("B140").Value = 855
$ws.Range("D140").Value = 797
$ws.Range("E140").Value = 6
$ws.Range("H140").Value = 52

# Row 212 (now "Seychelles" after the label swap above)
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0

# Row 213 (now "Montserrat" after the label swap above)
$ws.Range("D213").Value = 10
$ws.Range("H213").Value = 1

# --- Timestamp update ----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Junio de 2020 a las 20:59"
